$d = $word.ActiveDocument

# The "SURAT PERNYATAAN" heading contains a paragraph whose single run holds
# the "${no_surat}" placeholder. Prefix that placeholder with a bold
# "No. " label so the line reads "No. ${no_surat}".
$target = $d.Content
$found = $target.Find.Execute("`${no_surat}", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

if ($found) {
    $insertPos = $target.Start

    # Seed a one-character run with formatting that differs from the
    # "${no_surat}" run (Bold = False) so it stays a distinct run instead of
    # being swallowed into the placeholder's run, then swap in the real
    # "No. " label text (it inherits the seed run's formatting).
    $seed = $d.Range($insertPos, $insertPos)
    $seed.InsertBefore("X")
    $seed = $d.Range($insertPos, $insertPos + 1)
    $seed.Font.Bold = $false
    $seed.Text = "No. "

    # Match the label's bold weight to the rest of the line; the font,
    # size, language and no-proofing attributes are already inherited
    # correctly from the surrounding text.
    $label = $d.Range($insertPos, $insertPos + 4)
    $label.Font.Bold = $true
}
